$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates: preserve original style/text-type by
# capturing the cell's Style, forcing a Text number format just long
# enough to write the numeric-looking literal without Excel coercing
# it to a Number, then restoring the original Style object. ---
$cell = $ws.Range('D2')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '64.581.84'
$cell.Style = $origStyle

$cell = $ws.Range('D3')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.151.34'
$cell.Style = $origStyle

$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '613.09'
$cell.Style = $origStyle

$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '148.08'
$cell.Style = $origStyle

$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.150.13'
$cell.Style = $origStyle

$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.46'
$cell.Style = $origStyle

$cell = $ws.Range('D13')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0000258'
$cell.Style = $origStyle

$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '35.70'
$cell.Style = $origStyle

$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.674.00'
$cell.Style = $origStyle

$cell = $ws.Range('D17')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '64.523.00'
$cell.Style = $origStyle

$cell = $ws.Range('D18')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.156.28'
$cell.Style = $origStyle

$cell = $ws.Range('D19')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.93'
$cell.Style = $origStyle

$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '480.99'
$cell.Style = $origStyle

$cell = $ws.Range('D21')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '14.69'
$cell.Style = $origStyle

$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.02'
$cell.Style = $origStyle

$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '13.73'
$cell.Style = $origStyle

$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '84.04'
$cell.Style = $origStyle

$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.55'
$cell.Style = $origStyle

$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.120'
$cell.Style = $origStyle

$cell = $ws.Range('D34')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '26.44'
$cell.Style = $origStyle

$cell = $ws.Range('D38')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.21'
$cell.Style = $origStyle

$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '459.40'
$cell.Style = $origStyle

$cell = $ws.Range('D41')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0403'
$cell.Style = $origStyle

$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.121'
$cell.Style = $origStyle

$cell = $ws.Range('D43')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.40'
$cell.Style = $origStyle

$cell = $ws.Range('D44')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.848.68'
$cell.Style = $origStyle

$cell = $ws.Range('D47')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.46'
$cell.Style = $origStyle

$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '26.58'
$cell.Style = $origStyle

$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '120.33'
$cell.Style = $origStyle

# --- Volume(1h) (column E) updates: plain text assignment (values
# always contain '%' and padding spaces so Excel keeps them as text). ---
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('E3').Value = '  -1.12%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('E5').Value = '  +1.61%  '
$ws.Range('E6').Value = '  -3.31%  '
$ws.Range('E8').Value = '  -1.11%  '
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('E10').Value = '  -2.09%  '
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('E13').Value = '  -1.70%  '
$ws.Range('E14').Value = '  -4.39%  '
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('E19').Value = '  -2.57%  '
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('E21').Value = '  -1.58%  '
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('E23').Value = '  +2.09%  '
$ws.Range('E24').Value = '  -2.72%  '
$ws.Range('E25').Value = '  -1.28%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  -3.80%  '
$ws.Range('E28').Value = '  -3.22%  '
$ws.Range('E29').Value = '  -1.69%  '
$ws.Range('E30').Value = '  -2.48%  '
$ws.Range('E31').Value = '  -7.94%  '
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('E33').Value = '  -0.63%  '
$ws.Range('E34').Value = '  -2.37%  '
$ws.Range('E35').Value = '  +1.50%  '
$ws.Range('E36').Value = '  +6.09%  '
$ws.Range('E37').Value = '  -2.75%  '
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('E39').Value = '  -2.75%  '
$ws.Range('E40').Value = '  -1.33%  '
$ws.Range('E41').Value = '  -1.36%  '
$ws.Range('E42').Value = '  -6.01%  '
$ws.Range('E43').Value = '  -2.20%  '
$ws.Range('E44').Value = '  -2.98%  '
$ws.Range('E45').Value = '  -6.06%  '
$ws.Range('E46').Value = '  -3.64%  '
$ws.Range('E47').Value = '  +3.58%  '
$ws.Range('E48').Value = '  -3.58%  '
$ws.Range('E50').Value = '  -2.13%  '
$ws.Range('E51').Value = '  -0.44%  '
